$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29-87 down to 30-88.
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the weekly record.
$ws.Range("A29").Value = 9
$ws.Range("B29").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44536
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 100112022
$ws.Range("G29").Value = "Arveja Verde"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 43
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 16512
$ws.Range("N29").Value = "$/saco 25 kilos"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 660
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"

# Match the date-cell number format used by the rest of column D.
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
